# Boon Supply test-data refresh: rename the "Matthew *" test fundraiser to
# "Dwayne *" (first name + derived e-mails) on Scenario7, mirror the e-mail
# rename on the duplicate lookup sheets (Scenario8 / Scenario9), shift the
# StartMonth/StartDate/EndMonth/EndDate fundraiser windows on Scenario7, and
# update the active-sheet / selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Scenario7: FirstName + Email + fundraiser date-range updates
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Scenario7")

$ws7.Range("C2").Value = "Dwayne"
$ws7.Range("C3").Value = "Dwayne"
$ws7.Range("C4").Value = "Dwayne"
$ws7.Range("C5").Value = "Dwayne"

$ws7.Range("F2").Value = "dwaynewade@getnada.com"
$ws7.Range("F3").Value = "dwaynejohn@getnada.com"
$ws7.Range("F4").Value = "dwaynescott@getnada.com"
$ws7.Range("F5").Value = "dwaynejhonson@getnada.com"

# StartMonth / StartDate / EndMonth / EndDate - written column by column
# (top to bottom, then next column) so newly-introduced shared strings land
# in the same append order the source workbook ended up with.
$ws7.Range("O2").Value = "October"
$ws7.Range("O3").Value = "November"
$ws7.Range("O4").Value = "December"
$ws7.Range("O5").Value = "December"

$ws7.Range("P2").Value = "29"
$ws7.Range("P3").Value = "2"
$ws7.Range("P4").Value = "5"
$ws7.Range("P5").Value = "12"

$ws7.Range("Q2").Value = "November"
$ws7.Range("Q3").Value = "December"
$ws7.Range("Q4").Value = "April"
$ws7.Range("Q5").Value = "January"

$ws7.Range("R2").Value = "12"
$ws7.Range("R3").Value = "21"
$ws7.Range("R4").Value = "22"
$ws7.Range("R5").Value = "2"

# ---------------------------------------------------------------------------
# Scenario8 / Scenario9: same e-mail rename on their lookup/email column
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Scenario8")
$ws8.Range("C2").Value = "dwaynewade@getnada.com"
$ws8.Range("C3").Value = "dwaynejohn@getnada.com"
$ws8.Range("C4").Value = "dwaynescott@getnada.com"
$ws8.Range("C5").Value = "dwaynejhonson@getnada.com"

$ws9 = $wb.Worksheets.Item("Scenario9")
$ws9.Range("C2").Value = "dwaynewade@getnada.com"
$ws9.Range("C3").Value = "dwaynejohn@getnada.com"
$ws9.Range("C4").Value = "dwaynescott@getnada.com"
$ws9.Range("C5").Value = "dwaynejhonson@getnada.com"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping.
# Order matters: the last sheet whose Range(...).Select() runs becomes the
# workbook's active tab, so select Scenario9 and Scenario8 first, then
# finish on Scenario7 (which should end up tabSelected / activeTab).
# ---------------------------------------------------------------------------
$ws9.Range("C10").Select()
$ws8.Range("C2:D5").Select()
$ws7.Range("F16").Select()
